# Hospital Management System - DataDictionary.xlsx
# Adds the "Department" and "Doctor" table definitions to the data
# dictionary worksheet (employee repository / enum class support).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Department table (rows 61-65)
# ---------------------------------------------------------------------
$ws.Range("A61").Value = "Department"
$ws.Range("B61").Value = "departmentid"
$ws.Range("C61").Value = "int"
$ws.Range("E61").Value = "no"
$ws.Range("F61").Value = "pk"

$ws.Range("A62").Value = "Department"
$ws.Range("B62").Value = "hospitalid"
$ws.Range("C62").Value = "int"
$ws.Range("E62").Value = "no"
$ws.Range("F62").Value = "fk_hospital"

$ws.Range("A63").Value = "Department"
$ws.Range("B63").Value = "departmentcode"
$ws.Range("C63").Value = "int"
$ws.Range("E63").Value = "no"

$ws.Range("A64").Value = "Department"
$ws.Range("B64").Value = "departmentname"
$ws.Range("C64").Value = "varchar"
$ws.Range("D64").Value = 255
$ws.Range("E64").Value = "no"

$ws.Range("A65").Value = "Department"
$ws.Range("B65").Value = "discription"
$ws.Range("C65").Value = "varchar"
$ws.Range("D65").Value = "max"
$ws.Range("E65").Value = "yes"

# ---------------------------------------------------------------------
# Doctor table (rows 71-81)
# ---------------------------------------------------------------------
$ws.Range("A71").Value = "Doctor"
$ws.Range("B71").Value = "empid"
$ws.Range("C71").Value = "int "
$ws.Range("E71").Value = "no "
$ws.Range("F71").Value = "pk"

$ws.Range("A72").Value = "Doctor"
$ws.Range("B72").Value = "hospitalid"
$ws.Range("C72").Value = "int"
$ws.Range("E72").Value = "no"
$ws.Range("F72").Value = "fk_hospital"

$ws.Range("A73").Value = "Doctor"
$ws.Range("B73").Value = "departmentid"
$ws.Range("C73").Value = "int"
$ws.Range("E73").Value = "no"
$ws.Range("F73").Value = "fk_dipartment"

$ws.Range("A74").Value = "Doctor"
$ws.Range("B74").Value = "emptype"
$ws.Range("C74").Value = "int"
$ws.Range("E74").Value = "no"

$ws.Range("A75").Value = "Doctor"
$ws.Range("B75").Value = "empcode"
$ws.Range("C75").Value = "int"
$ws.Range("D75").Value = 20
$ws.Range("E75").Value = "no"

$ws.Range("A76").Value = "Doctor"
$ws.Range("B76").Value = "empname"
$ws.Range("C76").Value = "varchar"
$ws.Range("D76").Value = 255
$ws.Range("E76").Value = "no"

$ws.Range("A77").Value = "Doctor"
$ws.Range("B77").Value = "gender"
$ws.Range("C77").Value = "varchar"
$ws.Range("E77").Value = "no"

$ws.Range("A78").Value = "Doctor"
$ws.Range("B78").Value = "mobileno"
$ws.Range("C78").Value = "varchar"
$ws.Range("D78").Value = 15
$ws.Range("E78").Value = "no"

$ws.Range("A79").Value = "Doctor"
$ws.Range("B79").Value = "email"
$ws.Range("C79").Value = "varchar"
$ws.Range("D79").Value = 50
$ws.Range("E79").Value = "no"

$ws.Range("A80").Value = "Doctor"
$ws.Range("B80").Value = "qualification"
$ws.Range("C80").Value = "varchar"
$ws.Range("D80").Value = 255
$ws.Range("E80").Value = "no "

$ws.Range("A81").Value = "Doctor"
$ws.Range("B81").Value = "job specification"
$ws.Range("C81").Value = "varchar"
$ws.Range("D81").Value = 255
$ws.Range("E81").Value = "yes"

# ---------------------------------------------------------------------
# View state: the last-edited cell ends up selected, with the sheet
# scrolled so the new rows are visible.
# ---------------------------------------------------------------------
$ws.Range("B81").Select()
